# "updated docs for pages" - CodeSystem-nmdp-language-cs.xlsx
#
# Metadata sheet changes:
#   - B7  (Experimental)    : "false" -> "true"
#   - B8  (Date)             : "2023-01-12T09:36:27-06:00" -> "2023-02-16T14:43:10-06:00"
#   - B14 (Case Sensitive)   : (blank) -> "true"
#
# Values are written with a leading apostrophe so Excel stores them as text
# (shared strings) instead of auto-coercing "true"/"false" into native
# booleans. The formatting is then restored from an untouched sibling cell
# (B13) via Copy/PasteSpecial-formats so the original cell style (s="2") is
# preserved rather than picking up a stray "quote prefix" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

function Set-TextValue($cell, $text) {
    $target = $ws.Range($cell)
    $target.Value = "'" + $text
    $ws.Range("B13").Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

Set-TextValue "B7" "true"
Set-TextValue "B8" "2023-02-16T14:43:10-06:00"
Set-TextValue "B14" "true"

$excel.CutCopyMode = 0
